# Auto-generated Excel COM-interop script applying the Leviathan_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 7068.5
$ws.Range("I18").Value = 7068.5
$ws.Range("K18").Value = 7068.5
$ws.Range("M18").Value = -6784.5
$ws.Range("H19").Value = 5631.778
$ws.Range("I19").Value = 6210.75
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 6210.75
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = -6035.75
$ws.Range("N19").Value = -1350
$ws.Range("H43").Value = 7143.8887
$ws.Range("I43").Value = 3748.75
$ws.Range("J43").Value = 9860
$ws.Range("K43").Value = 3748.75
$ws.Range("L43").Value = 9860
$ws.Range("M43").Value = -3679.75
$ws.Range("N43").Value = -9998
$ws.Range("H88").Value = 2818.0908
$ws.Range("I88").Value = 3333.3333
$ws.Range("J88").Value = 2624.875
$ws.Range("K88").Value = 3333.3333
$ws.Range("L88").Value = 2624.875
$ws.Range("M88").Value = -2927.3333
$ws.Range("N88").Value = -3436.875
$ws.Range("H91").Value = 2818.0908
$ws.Range("I91").Value = 3333.3333
$ws.Range("J91").Value = 2624.875
$ws.Range("K91").Value = 3333.3333
$ws.Range("L91").Value = 2624.875
$ws.Range("M91").Value = -1929.3333
$ws.Range("N91").Value = -5432.875
$ws.Range("H107").Value = 31854.625
$ws.Range("I107").Value = 603.3
$ws.Range("J107").Value = 188111.25
$ws.Range("K107").Value = 603.3
$ws.Range("L107").Value = 188111.25
$ws.Range("M107").Value = 1316.7
$ws.Range("N107").Value = -191951.25
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H137").Value = 1696.3529
$ws.Range("I137").Value = 1794.7826
$ws.Range("J137").Value = 1490.5454
$ws.Range("K137").Value = 5384.3478
$ws.Range("L137").Value = 4471.6362
$ws.Range("M137").Value = -2834.3478
$ws.Range("N137").Value = -9571.6362
$ws.Range("H141").Value = 24518.738
$ws.Range("I141").Value = 24882.2
$ws.Range("K141").Value = 74646.6
$ws.Range("M141").Value = -69466.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1993.9231
$ws.Range("I2").Value = 1902.381
$ws.Range("K2").Value = 1902.381
$ws.Range("M2").Value = -1789.381
$ws.Range("H32").Value = 32298.816
$ws.Range("I32").Value = 5952.476
$ws.Range("J32").Value = 190376.86
$ws.Range("K32").Value = 5952.476
$ws.Range("L32").Value = 190376.86
$ws.Range("M32").Value = -5665.476
$ws.Range("N32").Value = -190950.86
$ws.Range("H61").Value = 1722.2307
$ws.Range("I61").Value = 1657.4166
$ws.Range("K61").Value = 1657.4166
$ws.Range("M61").Value = -1445.4166
$ws.Range("H86").Value = 30314
$ws.Range("J86").Value = 30314
$ws.Range("L86").Value = 30314
$ws.Range("N86").Value = -32686
$ws.Range("H89").Value = 30314
$ws.Range("J89").Value = 30314
$ws.Range("L89").Value = 90942
$ws.Range("N89").Value = -102798
$ws.Range("H116").Value = 1993.9231
$ws.Range("I116").Value = 1902.381
$ws.Range("K116").Value = 1902.381
$ws.Range("M116").Value = 391.6189999999999
$ws.Range("H136").Value = 1722.2307
$ws.Range("I136").Value = 1657.4166
$ws.Range("K136").Value = 4972.2498
$ws.Range("M136").Value = -2422.2498

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1993.9231
$ws.Range("I3").Value = 1902.381
$ws.Range("K3").Value = 1902.381
$ws.Range("M3").Value = -1788.381
$ws.Range("H134").Value = 1763.3077
$ws.Range("I134").Value = 1393.4
$ws.Range("K134").Value = 4180.200000000001
$ws.Range("M134").Value = -1645.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 23333668
$ws.Range("I3").Value = 23333668
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 23333668
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -23333555
$ws.Range("N3").ClearContents()
$ws.Range("H31").Value = 41831.35
$ws.Range("I31").Value = 40971.27
$ws.Range("K31").Value = 40971.27
$ws.Range("M31").Value = -40676.27
$ws.Range("H34").Value = 41831.35
$ws.Range("I34").Value = 40971.27
$ws.Range("K34").Value = 40971.27
$ws.Range("M34").Value = -40769.27
$ws.Range("H58").Value = 1863.3077
$ws.Range("I58").Value = 831
$ws.Range("J58").Value = 2748.1428
$ws.Range("K58").Value = 831
$ws.Range("L58").Value = 2748.1428
$ws.Range("M58").Value = -628
$ws.Range("N58").Value = -3154.1428
$ws.Range("H86").Value = 5999.3335
$ws.Range("H89").Value = 5999.3335
$ws.Range("H107").Value = 1755.6316
$ws.Range("I107").Value = 1698.091
$ws.Range("J107").Value = 1834.75
$ws.Range("K107").Value = 1698.091
$ws.Range("L107").Value = 1834.75
$ws.Range("M107").Value = 221.9090000000001
$ws.Range("N107").Value = -5674.75
$ws.Range("H136").Value = 1863.3077
$ws.Range("I136").Value = 831
$ws.Range("J136").Value = 2748.1428
$ws.Range("K136").Value = 2493
$ws.Range("L136").Value = 8244.4284
$ws.Range("M136").Value = 57
$ws.Range("N136").Value = -13344.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1830.8823
$ws.Range("J132").Value = 1860.3846
$ws.Range("L132").Value = 16743.4614
$ws.Range("N132").Value = -21803.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 16229.5
$ws.Range("H126").Value = 3475.111
$ws.Range("I126").Value = 2999.5
$ws.Range("J126").Value = 3855.6
$ws.Range("K126").Value = 8998.5
$ws.Range("L126").Value = 11566.8
$ws.Range("M126").Value = -6528.5
$ws.Range("N126").Value = -16506.8
$ws.Range("H132").Value = 1896.1333
$ws.Range("I132").Value = 1926.6923
$ws.Range("K132").Value = 5780.0769
$ws.Range("M132").Value = -3250.0769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 5093.129
$ws.Range("I20").Value = 4192.5713
$ws.Range("K20").Value = 4192.5713
$ws.Range("M20").Value = -3966.5713
$ws.Range("H100").Value = 41331.445
$ws.Range("I100").Value = 2716.6667
$ws.Range("J100").Value = 60638.832
$ws.Range("K100").Value = 2716.6667
$ws.Range("L100").Value = 60638.832
$ws.Range("M100").Value = -2175.6667
$ws.Range("N100").Value = -61720.832
$ws.Range("H122").Value = 3419.24
$ws.Range("I122").Value = 2676.0667
$ws.Range("K122").Value = 8028.2001
$ws.Range("M122").Value = -5578.2001
$ws.Range("H136").Value = 3438.3142
$ws.Range("I136").Value = 3064
$ws.Range("J136").Value = 5247.5
$ws.Range("K136").Value = 9192
$ws.Range("L136").Value = 15742.5
$ws.Range("M136").Value = -6642
$ws.Range("N136").Value = -20842.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 3740
$ws.Range("I6").Value = 10000
$ws.Range("J6").Value = 2175
$ws.Range("K6").Value = 10000
$ws.Range("L6").Value = 2175
$ws.Range("M6").Value = -9885
$ws.Range("N6").Value = -2405
$ws.Range("H21").Value = 7240
$ws.Range("I21").Value = 1200
$ws.Range("J21").Value = 8750
$ws.Range("K21").Value = 1200
$ws.Range("L21").Value = 8750
$ws.Range("M21").Value = -965
$ws.Range("N21").Value = -9220
$ws.Range("H35").Value = 7240
$ws.Range("I35").Value = 1200
$ws.Range("J35").Value = 8750
$ws.Range("K35").Value = 1200
$ws.Range("L35").Value = 8750
$ws.Range("M35").Value = -910
$ws.Range("N35").Value = -9330
$ws.Range("H62").Value = 37644.76
$ws.Range("I62").Value = 11378.571
$ws.Range("K62").Value = 11378.571
$ws.Range("M62").Value = -10754.571
$ws.Range("H65").Value = 37644.76
$ws.Range("I65").Value = 11378.571
$ws.Range("K65").Value = 56892.855
$ws.Range("M65").Value = -53772.855
$ws.Range("H122").Value = 2149.2222
$ws.Range("I122").Value = 1398.8334
$ws.Range("K122").Value = 4196.5002
$ws.Range("M122").Value = -1746.5002

Write-Host "Applied all Leviathan_Profits.xlsx updates"
